# Final Code Update Checkin For 207 Project
# Removes the recall / history-plot / confusion-matrix columns (P:U),
# repurposes the remaining N/O header+data columns, and refreshes the
# accuracy/model-size/history-plot values for each run row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing columns (Test Accuracy was column N and is
# being removed; Train/Val Recall, Drowsy/Non-Drowsy Recall, History Plot
# and Confusion Matrix - columns P through U - go away entirely). This
# shifts the remaining data left and shrinks the used range to A1:O9.
$ws.Range("P1:U9").Delete()

# Re-label the remaining two trailing headers.
$ws.Range("N1").Value = "Model Size"
$ws.Range("O1").Value = "History Plot"

# Refresh per-run data: Train Accuracy (L), Validation Accuracy (M),
# Model Size (N, now the old O-column value) and History Plot path (O,
# re-indexed to start at 0).
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 189.1264686584473
$ws.Range("O2").Value = "./plots/base_full_face/hist/history_0.png"

$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 189.1264686584473
$ws.Range("O3").Value = "./plots/base_full_face/hist/history_1.png"

$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.949999988079071
$ws.Range("N4").Value = 126.5555458068848
$ws.Range("O4").Value = "./plots/base_full_face/hist/history_2.png"

$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9833333492279053
$ws.Range("N5").Value = 126.5555458068848
$ws.Range("O5").Value = "./plots/base_full_face/hist/history_3.png"

$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9833333492279053
$ws.Range("N6").Value = 151.3449745178223
$ws.Range("O6").Value = "./plots/base_full_face/hist/history_4.png"

$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 151.3449745178223
$ws.Range("O7").Value = "./plots/base_full_face/hist/history_5.png"

$ws.Range("L8").Value = 0.9937499761581421
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 101.2740516662598
$ws.Range("O8").Value = "./plots/base_full_face/hist/history_6.png"

$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 101.2740516662598
$ws.Range("O9").Value = "./plots/base_full_face/hist/history_7.png"
